$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 42200
$ws.Range("J3").Value = 42200
$ws.Range("L3").Value = 42200
$ws.Range("N3").Value = -42428
$ws.Range("H7").Value = 8006
$ws.Range("J7").Value = 8006
$ws.Range("L7").Value = 8006
$ws.Range("N7").Value = -8230
$ws.Range("H13").Value = 5006
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 5006
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 5006
$ws.Range("M13").ClearContents()
$ws.Range("N13").Value = -5344
$ws.Range("H14").Value = 8006
$ws.Range("J14").Value = 8006
$ws.Range("L14").Value = 8006
$ws.Range("N14").Value = -8388
$ws.Range("H40").Value = 1690
$ws.Range("I40").Value = 1475
$ws.Range("J40").Value = 1833.3334
$ws.Range("K40").Value = 1475
$ws.Range("L40").Value = 1833.3334
$ws.Range("M40").Value = -1300
$ws.Range("N40").Value = -2183.3334
$ws.Range("H44").Value = 18000
$ws.Range("J44").Value = 18000
$ws.Range("L44").Value = 18000
$ws.Range("N44").Value = -18924
$ws.Range("H102").Value = 42200
$ws.Range("J102").Value = 42200
$ws.Range("L102").Value = 42200
$ws.Range("N102").Value = -48690
$ws.Range("H111").Value = 1100
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H118").Value = 2420.4285
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 2420.4285
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 7261.2855
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -10575.2855
$ws.Range("H129").Value = 1010.2059
$ws.Range("J129").Value = 1226.5834
$ws.Range("L129").Value = 3679.7502
$ws.Range("N129").Value = -13679.7502
$ws.Range("H132").Value = 2092.38
$ws.Range("I132").Value = 1531.2632
$ws.Range("K132").Value = 4593.7896
$ws.Range("M132").Value = -2063.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 700.54346
$ws.Range("I2").Value = 403.05713
$ws.Range("K2").Value = 403.05713
$ws.Range("M2").Value = -290.05713
$ws.Range("H32").Value = 16207.51
$ws.Range("I32").Value = 16437.459
$ws.Range("J32").Value = 14000
$ws.Range("K32").Value = 16437.459
$ws.Range("L32").Value = 14000
$ws.Range("M32").Value = -16150.459
$ws.Range("N32").Value = -14574
$ws.Range("H116").Value = 700.54346
$ws.Range("I116").Value = 403.05713
$ws.Range("K116").Value = 403.05713
$ws.Range("M116").Value = 1890.94287
$ws.Range("H122").Value = 1729.4584
$ws.Range("I122").Value = 1822.8422
$ws.Range("J122").Value = 1374.6
$ws.Range("K122").Value = 5468.5266
$ws.Range("L122").Value = 4123.799999999999
$ws.Range("M122").Value = -3018.5266
$ws.Range("N122").Value = -9023.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 700.54346
$ws.Range("I3").Value = 403.05713
$ws.Range("K3").Value = 403.05713
$ws.Range("M3").Value = -289.05713
$ws.Range("H92").Value = 199067.67
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 199067.67
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 199067.67
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -204059.67
$ws.Range("H94").Value = 805.3333
$ws.Range("I94").Value = 637.1667
$ws.Range("K94").Value = 637.1667
$ws.Range("M94").Value = -186.1667
$ws.Range("H95").Value = 44333.332
$ws.Range("J95").Value = 44333.332
$ws.Range("L95").Value = 44333.332
$ws.Range("N95").Value = -49825.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27030002
$ws.Range("I31").Value = 37039280
$ws.Range("K31").Value = 37039280
$ws.Range("M31").Value = -37038985
$ws.Range("H32").Value = 7007.3335
$ws.Range("J32").Value = 10011
$ws.Range("L32").Value = 10011
$ws.Range("N32").Value = -10643
$ws.Range("H34").Value = 27030002
$ws.Range("I34").Value = 37039280
$ws.Range("K34").Value = 37039280
$ws.Range("M34").Value = -37039078
$ws.Range("H97").Value = 39193.5
$ws.Range("J97").Value = 39193.5
$ws.Range("L97").Value = 39193.5
$ws.Range("N97").Value = -41175.5
$ws.Range("H99").Value = 2202.1904
$ws.Range("I99").Value = 2378.9412
$ws.Range("J99").Value = 1451
$ws.Range("K99").Value = 2378.9412
$ws.Range("L99").Value = 1451
$ws.Range("M99").Value = -880.9412000000002
$ws.Range("N99").Value = -4447
$ws.Range("H122").Value = 992.8461
$ws.Range("I122").Value = 935.9091
$ws.Range("J122").Value = 1306
$ws.Range("K122").Value = 2807.7273
$ws.Range("L122").Value = 3918
$ws.Range("M122").Value = -357.7273
$ws.Range("N122").Value = -8818
$ws.Range("H126").Value = 2202.1904
$ws.Range("I126").Value = 2378.9412
$ws.Range("J126").Value = 1451
$ws.Range("K126").Value = 7136.823600000001
$ws.Range("L126").Value = 4353
$ws.Range("M126").Value = -4666.823600000001
$ws.Range("N126").Value = -9293
$ws.Range("H132").Value = 2720.65
$ws.Range("I132").Value = 2275.875
$ws.Range("K132").Value = 6827.625
$ws.Range("M132").Value = -4297.625
$ws.Range("H134").Value = 1303.5
$ws.Range("I134").Value = 1275.96
$ws.Range("J134").Value = 1533
$ws.Range("K134").Value = 3827.88
$ws.Range("L134").Value = 4599
$ws.Range("M134").Value = -1292.88
$ws.Range("N134").Value = -9669

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 873.79
$ws.Range("I131").Value = 503.33334
$ws.Range("J131").Value = 885.24744
$ws.Range("K131").Value = 1510.00002
$ws.Range("L131").Value = 2655.74232
$ws.Range("M131").Value = 3529.99998
$ws.Range("N131").Value = -12735.74232
$ws.Range("H134").Value = 4527.5815
$ws.Range("I134").Value = 1704.9445
$ws.Range("J134").Value = 6559.88
$ws.Range("K134").Value = 5114.833500000001
$ws.Range("L134").Value = 19679.64
$ws.Range("M134").Value = -44.83350000000064
$ws.Range("N134").Value = -29819.64

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3634.2273
$ws.Range("I122").Value = 3786.8333
$ws.Range("J122").Value = 2947.5
$ws.Range("K122").Value = 11360.4999
$ws.Range("L122").Value = 8842.5
$ws.Range("M122").Value = -8910.499899999999
$ws.Range("N122").Value = -13742.5
$ws.Range("H123").Value = 11290.25
$ws.Range("J123").Value = 11290.25
$ws.Range("L123").Value = 11290.25
$ws.Range("N123").Value = -16190.25
$ws.Range("H126").Value = 2286.963
$ws.Range("I126").Value = 2118.7896
$ws.Range("J126").Value = 2686.375
$ws.Range("K126").Value = 6356.3688
$ws.Range("L126").Value = 8059.125
$ws.Range("M126").Value = -3886.3688
$ws.Range("N126").Value = -12999.125
$ws.Range("H132").Value = 2937.2273
$ws.Range("I132").Value = 1973.8667
$ws.Range("K132").Value = 5921.6001
$ws.Range("M132").Value = -3391.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 13316.263
$ws.Range("I61").Value = 18238.846
$ws.Range("J61").Value = 2650.6667
$ws.Range("K61").Value = 18238.846
$ws.Range("L61").Value = 2650.6667
$ws.Range("M61").Value = -18036.846
$ws.Range("N61").Value = -3054.6667
$ws.Range("H113").Value = 13316.263
$ws.Range("I113").Value = 18238.846
$ws.Range("J113").Value = 2650.6667
$ws.Range("K113").Value = 18238.846
$ws.Range("L113").Value = 2650.6667
$ws.Range("M113").Value = -16068.846
$ws.Range("N113").Value = -6990.6667
$ws.Range("H122").Value = 12505249
$ws.Range("I122").Value = 20839708
$ws.Range("J122").Value = 3561.875
$ws.Range("K122").Value = 62519124
$ws.Range("L122").Value = 10685.625
$ws.Range("M122").Value = -62516674
$ws.Range("N122").Value = -15585.625
$ws.Range("H132").Value = 6027.115
$ws.Range("I132").Value = 6073.913
$ws.Range("J132").Value = 5668.3335
$ws.Range("K132").Value = 18221.739
$ws.Range("L132").Value = 17005.0005
$ws.Range("M132").Value = -15691.739
$ws.Range("N132").Value = -22065.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3033.8
$ws.Range("I62").Value = 2917.25
$ws.Range("J62").Value = 3500
$ws.Range("K62").Value = 2917.25
$ws.Range("L62").Value = 3500
$ws.Range("M62").Value = -2293.25
$ws.Range("N62").Value = -4748
$ws.Range("H65").Value = 3033.8
$ws.Range("I65").Value = 2917.25
$ws.Range("J65").Value = 3500
$ws.Range("K65").Value = 14586.25
$ws.Range("L65").Value = 17500
$ws.Range("M65").Value = -11466.25
$ws.Range("N65").Value = -23740
$ws.Range("H75").Value = 73333.336
$ws.Range("I75").Value = 200000
$ws.Range("J75").Value = 10000
$ws.Range("K75").Value = 200000
$ws.Range("L75").Value = 10000
$ws.Range("M75").Value = -199064
$ws.Range("N75").Value = -11872
$ws.Range("H78").Value = 73333.336
$ws.Range("I78").Value = 200000
$ws.Range("J78").Value = 10000
$ws.Range("K78").Value = 600000
$ws.Range("L78").Value = 30000
$ws.Range("M78").Value = -595320
$ws.Range("N78").Value = -39360
$ws.Range("H113").Value = 684
$ws.Range("I113").Value = 426.9524
$ws.Range("K113").Value = 1280.8572
$ws.Range("M113").Value = 889.1428000000001
$ws.Range("H114").Value = 80000
$ws.Range("J114").Value = 80000
$ws.Range("L114").Value = 80000
$ws.Range("N114").Value = -88678
$ws.Range("H122").Value = 35716988
$ws.Range("I122").Value = 62500700
$ws.Range("J122").Value = 5370
$ws.Range("K122").Value = 187502100
$ws.Range("L122").Value = 16110
$ws.Range("M122").Value = -187499650
$ws.Range("N122").Value = -21010
$ws.Range("H132").Value = 2809.5
$ws.Range("I132").Value = 2809.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 8428.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5898.5
$ws.Range("N132").ClearContents()
$ws.Range("H141").Value = 70168.60000000001
$ws.Range("J141").Value = 70168.60000000001
$ws.Range("L141").Value = 70168.60000000001
$ws.Range("N141").Value = -80528.60000000001
